# MOH 515 (Post Outbreak) — "worked on household indicators"
#
# 1. Re-label the existing "Form Summary" group's question labels to
#    Title Case, rename the group's own name from "inputs" to
#    "form_summary", and reword the CHP-reported-total label.
# 2. Append a brand-new "household_indicators" group (rows 10-16) with
#    five new integer questions.
# 3. Widen column C (label) and give column B (name) an explicit width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- 1. Re-label existing "Form Summary" group -----------------------
$ws.Range("B2").Value = "form_summary"
$ws.Range("C3").Value = "What Is Your Name?"
$ws.Range("C4").Value = "What Is Your Area?"
$ws.Range("C5").Value = "What Is Your Linked Facility?"
$ws.Range("C6").Value = "What County Do You Belong To?"
$ws.Range("C7").Value = "How Many CHPs Are In Your Area?"
$ws.Range("C8").Value = "How Many CHPs Submitted Monthly Report?"

# --- 2. Append the new "household_indicators" group ------------------

# Row 10: begin group / household_indicators / Household Indicators
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = "begin group"
$ws.Range("B10").Value = "household_indicators"
$ws.Range("C10").Value = "Household Indicators"

# Rows 11-13: integer questions, required=yes, appearence=numbers
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A11:D11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:D13").PasteSpecial(-4122) | Out-Null

$ws.Range("F7").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null

$ws.Range("A11").Value = "integer"
$ws.Range("B11").Value = "total_households"
$ws.Range("C11").Value = "Total Households In The Area?"
$ws.Range("D11").Value = "yes"
$ws.Range("F11").Value = "numbers"

$ws.Range("A12").Value = "integer"
$ws.Range("B12").Value = "new_households"
$ws.Range("C12").Value = "Number Of New Households Registered This Month?"
$ws.Range("D12").Value = "yes"
$ws.Range("F12").Value = "numbers"

$ws.Range("A13").Value = "integer"
$ws.Range("B13").Value = "new_households_visited"
$ws.Range("C13").Value = "Number Of New Households Visited This Month?"
$ws.Range("D13").Value = "yes"
$ws.Range("F13").Value = "numbers"

# Rows 14-15: integer questions, required=yes, no appearence
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A14:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("A15:D15").PasteSpecial(-4122) | Out-Null

$ws.Range("A14").Value = "integer"
$ws.Range("B14").Value = "new_households_with_clean_water"
$ws.Range("C14").Value = "Number Of New Households Visited This Month With Clean Water Access?"
$ws.Range("D14").Value = "yes"

$ws.Range("A15").Value = "integer"
$ws.Range("B15").Value = "new_households_with_latrines"
$ws.Range("D15").Value = "yes"

# C15's label keeps the odd "filled" style (s=2, same as H8) seen in the source
$ws.Range("H8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = "Number Of New Households Visited This Month With Latrines/Toilets?"

# Row 16: end group
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "end group"

# --- 3. Column widths --------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 59.3
